# Weekly fruit/hortaliza update:
# A new price observation (row) is inserted for "Vega Monumental Concepción -
# Zapallo italiano" ahead of the existing row 227, pushing the existing
# rows 227-251 down to 228-252. The new row carries the latest week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new blank row at position 227; everything currently at/after
# row 227 (through 251) shifts down one row (to 228-252), and the sheet's
# dimension grows to A1:R252 automatically.
$ws.Rows.Item(227).Insert()

# Populate the newly inserted row 227 with the new weekly record.
$ws.Range("A227").Value = 11
$ws.Range("B227").Value = "Vega Monumental Concepción"
$ws.Range("C227").Value = "Bíobío"
$ws.Range("D227").Value = 45132
$ws.Range("E227").Value = 8
$ws.Range("F227").Value = 100112032
$ws.Range("G227").Value = "Zapallo italiano"
$ws.Range("H227").Value = "Sin especificar"
$ws.Range("I227").Value = "Primera"
$ws.Range("J227").Value = 50
$ws.Range("K227").Value = 16000
$ws.Range("L227").Value = 16000
$ws.Range("M227").Value = 16000
$ws.Range("N227").Value = "$/caja 50 unidades"
$ws.Range("O227").Value = "Región de Arica y Parinacota"
$ws.Range("P227").Value = 320
$ws.Range("Q227").Value = 50
$ws.Range("R227").Value = "Hortaliza"
